$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '68.216.74'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.53%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.640.38'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.54%  '

# Row 4
$ws.Range("E4").Value = '  +0.03%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '597.99'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.12%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '154.56'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.67%  '

# Row 8
$ws.Range("E8").Value = '  -0.78%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.639.36'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.53%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.145'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +8.45%  '

# Row 11
$ws.Range("E11").Value = '  -0.67%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.24'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.46%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.353'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.68%  '

# Row 14
$ws.Range("B14").Value = 'Avalanche'
$ws.Range("C14").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.84'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.06%  '

# Row 15
$ws.Range("B15").Value = 'ShibaInu'
$ws.Range("C15").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000192'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.23%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.122.98'
$ws.Range("D16").Style = "Normal"

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '68.146.40'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.48%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.630.80'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.22%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.35'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.22%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '363.53'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.96%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.43'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.10%  '

# Row 22
$ws.Range("E22").Value = '  +2.75%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.83'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.85%  '

# Row 24
$ws.Range("E24").Value = '  -1.38%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '75.24'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +4.53%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.10%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.77'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.51%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0000105'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.54%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.777.07'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.60%  '

# Row 30
$ws.Range("E30").Value = '  -0.22%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '560.11'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.80%  '

# Row 32
$ws.Range("E32").Value = '  +1.60%  '

# Row 33
$ws.Range("E33").Value = '  -0.02%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.85'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.95%  '

# Row 35
$ws.Range("E35").Value = '  +1.52%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.00'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.03%  '

# Row 37
$ws.Range("E37").Value = '  +2.70%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '161.72'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.20%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '19.30'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.93%  '

# Row 40
$ws.Range("E40").Value = '  +1.52%  '

# Row 41
$ws.Range("E41").Value = '  -0.12%  '

# Row 42
$ws.Range("E42").Value = '  -0.41%  '

# Row 43
$ws.Range("E43").Value = '  +1.45%  '

# Row 44
$ws.Range("B44").Value = 'WhiteBITCoin'
$ws.Range("C44").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '17.75'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.27%  '

# Row 45
$ws.Range("B45").Value = 'dogwifhat'
$ws.Range("C45").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.63'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.05%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '40.62'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.31%  '

# Row 47
$ws.Range("E47").Value = '  -0.02%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '157.39'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.61%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.74'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.73%  '

# Row 50
$ws.Range("B50").Value = 'InjectiveProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '21.76'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.48%  '

# Row 51
$ws.Range("B51").Value = 'Optimism'
$ws.Range("C51").Value = 'https://coinranking.com/coin/n1p-s_gm1+optimism-op'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.69'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.00%  '
